# Weekly driver report update for 2025-04-21
# Update rows 4-6 of the "Bad Drivers" table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: now the "Intel(R) Dual Band Wireless-AC 8265 - 20.70.32.1" entry
$ws.Range("A4").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.32.1"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 151
$ws.Range("D4").Value = 98.40000000000001

# Row 5: now the "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4" entry
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B5").Value = 113
$ws.Range("C5").Value = 2780
$ws.Range("D5").Value = 98.8

# Row 6: updated totals
$ws.Range("B6").Value = 115
$ws.Range("C6").Value = 2933
